$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> source row number (data for columns D,J,K,L,M,O,P
# moves from the source row to the new row). Computed from the diff by
# matching the distinctive "Fecha" (D) serial values, and confirmed against
# every other changed column (J,K,L,M,O,P).
$map = @{
    2  = 19
    3  = 24
    4  = 12
    5  = 2
    6  = 4
    7  = 28
    8  = 17
    9  = 15
    10 = 21
    11 = 27
    12 = 13
    13 = 22
    14 = 10
    15 = 20
    16 = 25
    17 = 18
    18 = 29
    19 = 26
    20 = 8
    21 = 30
    22 = 23
    23 = 5
    24 = 7
    25 = 11
    26 = 6
    27 = 3
    28 = 14
    29 = 16
    30 = 9
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the columns that move, for every source
# row, before any writes happen (writes happen in-place on the same sheet).
# Value2 is used (rather than Value) to read/write the plain scalar.
$snapshot = @{}
foreach ($row in 2..30) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Apply the permutation: row N receives the snapshot of row $map[N].
foreach ($row in 2..30) {
    $srcRow = $map[$row]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
